$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.121.23'
$ws.Range('E2').Value = '  -3.92%  '
$ws.Range('D3').Value = '1.962.48'
$ws.Range('E3').Value = '  -4.06%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''241.47'
$ws.Range('E5').Value = '  -3.97%  '
$ws.Range('D6').Value = '''0.623'
$ws.Range('E6').Value = '  -3.94%  '
$ws.Range('D7').Value = '''59.95'
$ws.Range('E7').Value = '  -7.44%  '
$ws.Range('D9').Value = '''0.372'
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').Value = '''56.52'
$ws.Range('E10').Value = '  -4.68%  '
$ws.Range('D11').Value = '''0.0799'
$ws.Range('E11').Value = '  +5.88%  '
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('E13').Value = '  -5.91%  '
$ws.Range('D14').Value = '''22.21'
$ws.Range('E14').Value = '  +9.66%  '
$ws.Range('D15').Value = '''14.02'
$ws.Range('E15').Value = '  -7.02%  '
$ws.Range('D16').Value = '2.246.87'
$ws.Range('E16').Value = '  -4.22%  '
$ws.Range('D17').Value = '''5.39'
$ws.Range('E17').Value = '  -3.17%  '
$ws.Range('D18').Value = '1.971.31'
$ws.Range('E18').Value = '  -3.89%  '
$ws.Range('D19').Value = '36.001.73'
$ws.Range('E19').Value = '  -4.20%  '
$ws.Range('D20').Value = '''70.83'
$ws.Range('D21').Value = '0.0₃0855'
$ws.Range('E21').Value = '  -2.00%  '
$ws.Range('D22').Value = '''234.94'
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('E23').Value = '  -2.62%  '
$ws.Range('E24').Value = '  +0.30%  '
$ws.Range('E25').Value = '  -5.41%  '
$ws.Range('E26').Value = '  -4.00%  '
$ws.Range('D27').Value = '''9.75'
$ws.Range('E27').Value = '  +2.03%  '
$ws.Range('D28').Value = '''160.62'
$ws.Range('E28').Value = '  +0.89%  '
$ws.Range('D29').Value = '''19.76'
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('E30').Value = '  +14.13%  '
$ws.Range('E32').Value = '  -6.47%  '
$ws.Range('E33').Value = '  -5.31%  '
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('E35').Value = '  -7.31%  '
$ws.Range('D36').Value = '''6.25'
$ws.Range('E36').Value = '  +3.25%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  -6.59%  '
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('D40').Value = '''3.08'
$ws.Range('E40').Value = '  +10.25%  '
$ws.Range('D41').Value = '''0.0987'
$ws.Range('E41').Value = '  -4.38%  '
$ws.Range('E42').Value = '  -0.65%  '
$ws.Range('D43').Value = '''2.88'
$ws.Range('E43').Value = '  -2.04%  '
$ws.Range('E44').Value = '  -2.51%  '
$ws.Range('E45').Value = '  -4.56%  '
$ws.Range('D46').Value = '''92.09'
$ws.Range('E46').Value = '  -2.94%  '
$ws.Range('D47').Value = '''15.92'
$ws.Range('E47').Value = '  -5.43%  '
$ws.Range('E48').Value = '  -7.51%  '
$ws.Range('D49').Value = '1.333.31'
$ws.Range('E49').Value = '  -6.32%  '
$ws.Range('E50').Value = '  -3.93%  '
$ws.Range('D51').Value = '2.143.96'
$ws.Range('E51').Value = '  -3.92%  '
